$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its numeric-looking values stored as TEXT
# (matches the source data which stores prices as literal strings,
# some containing multiple "." thousand/decimal separators).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '25.769.43'
$ws.Range('E2').Value = '  -2.63%  '
$ws.Range('D3').Value = '1.744.08'
$ws.Range('E3').Value = '  -5.11%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '237.32'
$ws.Range('E5').Value = '  -9.24%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.5046'
$ws.Range('E7').Value = '  -6.35%  '
$ws.Range('D8').Value = '41.91'
$ws.Range('E8').Value = '  -6.63%  '
$ws.Range('D9').Value = '0.2661'
$ws.Range('E9').Value = '  -11.87%  '
$ws.Range('D10').Value = '0.06162'
$ws.Range('E10').Value = '  -10.30%  '
$ws.Range('D11').Value = '1.744.06'
$ws.Range('E11').Value = '  -5.05%  '
$ws.Range('D12').Value = '0.06921'
$ws.Range('E12').Value = '  -3.90%  '
$ws.Range('D13').Value = '15.45'
$ws.Range('E13').Value = '  -12.86%  '
$ws.Range('D14').Value = '4.510'
$ws.Range('E14').Value = '  -9.60%  '
$ws.Range('D15').Value = '0.5993'
$ws.Range('E15').Value = '  -18.84%  '
$ws.Range('D16').Value = '77.00'
$ws.Range('E16').Value = '  -13.76%  '
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = '25.774.60'
$ws.Range('E19').Value = '  -2.71%  '
$ws.Range('D20').Value = '0.000006817'
$ws.Range('E20').Value = '  -13.40%  '
$ws.Range('D21').Value = '11.57'
$ws.Range('E21').Value = '  -16.40%  '
$ws.Range('D22').Value = '1.966.57'
$ws.Range('E22').Value = '  -5.28%  '
$ws.Range('D23').Value = '4.056'
$ws.Range('E23').Value = '  -11.63%  '
$ws.Range('D24').Value = '5.207'
$ws.Range('E24').Value = '  -12.79%  '
$ws.Range('D25').Value = '8.146'
$ws.Range('D26').Value = '137.25'
$ws.Range('E26').Value = '  -3.83%  '
$ws.Range('D27').Value = '1.524'
$ws.Range('E27').Value = '  -9.79%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '15.01'
$ws.Range('E28').Value = '  -11.55%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '1.809'
$ws.Range('E29').Value = '  -17.54%  '
$ws.Range('D30').Value = '103.63'
$ws.Range('E30').Value = '  -6.24%  '
$ws.Range('D31').Value = '3.764'
$ws.Range('E31').Value = '  -10.91%  '
$ws.Range('D32').Value = '0.08115'
$ws.Range('E32').Value = '  -8.07%  '
$ws.Range('D33').Value = '3.471'
$ws.Range('E33').Value = '  -13.82%  '
$ws.Range('D34').Value = '0.04517'
$ws.Range('E34').Value = '  -6.11%  '
$ws.Range('D35').Value = '0.9998'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').Value = '2.649'
$ws.Range('E36').Value = '  -9.10%  '
$ws.Range('D37').Value = '0.9821'
$ws.Range('E37').Value = '  -13.13%  '
$ws.Range('D38').Value = '0.6109'
$ws.Range('E38').Value = '  -16.09%  '
$ws.Range('D39').Value = '2.670'
$ws.Range('E39').Value = '  -13.68%  '
$ws.Range('D40').Value = '0.01551'
$ws.Range('E40').Value = '  -9.34%  '
$ws.Range('D41').Value = '1.916'
$ws.Range('E41').Value = '  -15.16%  '
$ws.Range('D42').Value = '1.001'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').Value = '101.98'
$ws.Range('E43').Value = '  -5.32%  '
$ws.Range('D44').Value = '0.3816'
$ws.Range('E44').Value = '  -19.18%  '
$ws.Range('D45').Value = '5.079'
$ws.Range('E45').Value = '  -13.71%  '
$ws.Range('D46').Value = '0.7330'
$ws.Range('E46').Value = '  -18.94%  '
$ws.Range('D47').Value = '0.05363'
$ws.Range('E47').Value = '  -7.11%  '
$ws.Range('D48').Value = '0.1104'
$ws.Range('E48').Value = '  -10.98%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '30.17'
$ws.Range('E49').Value = '  -13.21%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').Value = '5.927'
$ws.Range('E50').Value = '  -19.75%  '
$ws.Range('D51').Value = '52.58'
$ws.Range('E51').Value = '  -12.23%  '
